# fix: standardize Notice u/s 94 line above To, across all templates in
# bold underline format.
#
# 1. Insert a new "Notice u/s 94 BNSS, 2023" paragraph (bold, underlined,
#    justified) directly above the "To," paragraph.
# 2. Remove the old, differently-formatted "Notice u/s 94 BNSS, 2023"
#    paragraph that used to sit directly above the "Subject: -" line.

$d = $word.ActiveDocument

# --- Step 1: insert the new paragraph right before the "To," paragraph ---

$toIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    $trimmed = $txt.Trim()
    if ($trimmed -eq "To,") {
        $toIndex = $i
        break
    }
}

if ($toIndex -eq 0) {
    throw "Could not locate the 'To,' paragraph"
}

$beforePara = $d.Paragraphs.Item($toIndex - 1)
$beforePara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($toIndex)
$newRange = $newPara.Range

$noticeXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
    '<w:p>' +
    '<w:pPr><w:jc w:val="both"/></w:pPr>' +
    '<w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>Notice u/s 94 BNSS, 2023</w:t></w:r>' +
    '</w:p>' +
    '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$newRange.InsertXML($noticeXml) | Out-Null

# --- Step 2: delete the old "Notice u/s 94 BNSS, 2023" paragraph that used
#     to live right above "Subject: -" ---

$oldIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $trimmed = $p.Range.Text.Trim()
    if ($trimmed -eq "Notice u/s 94 BNSS, 2023") {
        $nextTxt = $d.Paragraphs.Item($i + 1).Range.Text.Trim()
        if ($nextTxt.StartsWith("Subject")) {
            $oldIndex = $i
            break
        }
    }
}

if ($oldIndex -eq 0) {
    throw "Could not locate the legacy 'Notice u/s 94 BNSS, 2023' paragraph"
}

$d.Paragraphs.Item($oldIndex).Range.Delete()

Write-Host "Done: inserted new Notice paragraph above 'To,' and removed legacy Notice paragraph above 'Subject:'"
